$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths ---
$ws.Columns.Item(3).ColumnWidth = 18.7109375

# --- Row 9 header formatting: center align A9:C9 ---
$ws.Range("A9:C9").HorizontalAlignment = -4108   # xlCenter

# --- Row 13 ---
$ws.Range("A13").Value = "Num"                    # new shared string: Num
$ws.Range("C13").Formula = "=B5/A10"

# --- Row 15 label (typed alongside Num) ---
$ws.Range("A15").Value = "Deno"                   # new shared string: Deno
$ws.Range("B15").Value = 1
$ws.Range("C15").Formula = "=B10/A10"
$ws.Range("D15").Formula = "=C10/A10"

# --- Row 9 header text fix + row 14 labels (reuse the new S^0 string) ---
$ws.Range("C9").Value = "S^0"                     # new shared string: S^0
$ws.Range("B14").Value = "S^2"
$ws.Range("C14").Value = "S"
$ws.Range("D14").Value = "S^0"

# --- Column F labels (bold) + column G formulas ---
$ws.Range("F13").Value = "K"                      # new shared string: K
$ws.Range("F13").Font.Bold = $true
$ws.Range("G13").Formula = "=C13/D15"

$ws.Range("F14").Value = "Wn"                     # new shared string: Wn
$ws.Range("F14").Font.Bold = $true
$ws.Range("G14").Formula = "=SQRT(D15)"

$ws.Range("F15").Value = "Xi"                     # new shared string: Xi
$ws.Range("F15").Font.Bold = $true
$ws.Range("G15").Formula = "=C15/(2*G14)"

$ws.Range("F16").Value = "Ts"                     # new shared string: Ts
$ws.Range("F16").Font.Bold = $true
$ws.Range("G16").Formula = "=2*G15/G14"

$ws.Range("F17").Value = "Te"                     # new shared string: Te
$ws.Range("F17").Font.Bold = $true
$ws.Range("G17").Formula = "=G16*4"

# --- Row 18 ---
$ws.Range("A18").Value = "Regla de tres"          # new shared string: Regla de tres

# --- Row 19 ---
$ws.Range("A19").Value = "value"                  # new shared string: value
$ws.Range("B19").Value = "%"                      # new shared string: %
$ws.Range("E19").Font.Underline = $true

# --- Rows 20-23: regla de tres computations with a shared fill-down formula ---
$ws.Range("A20").Value = 22.06
$ws.Range("B20").Formula = "=A20*100/22.65"
$ws.Range("A21").Value = 14.226
$ws.Range("A22").Value = 14.82
$ws.Range("B21:B23").Formula = "=A21*100/22.65"

# --- View: zoom + selection ---
$excel.ActiveWindow.Zoom = 130
$ws.Range("G19").Select()
